$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new part entry
$ws.Range("A3").Value = "too "
$ws.Range("B3").Value = "cool"
$ws.Range("C3").Value = "4987321"
$ws.Range("D3").Value = "Mezzanine"
$ws.Range("E3").Value = "EN-48"
$ws.Range("F3").Value = "45"

# Row 4: another new part entry
$ws.Range("A4").Value = "a"
$ws.Range("B4").Value = "asf"
$ws.Range("C4").Value = "qwe"
$ws.Range("D4").Value = "PLC Room"
$ws.Range("E4").Value = "a-q"
$ws.Range("F4").Value = "qq"
